$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column style (format) from A301 down through A302:A328
$ws.Range("A301").Copy() | Out-Null
$ws.Range("A302:A328").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the data for rows 302-328 (update through 28 luglio)
$ws.Cells.Item(302, 1).Value2 = 44376
$ws.Cells.Item(302, 2).Value2 = 0
$ws.Cells.Item(302, 3).Value2 = 2
$ws.Cells.Item(302, 4).Value2 = 6.063729800200103
$ws.Cells.Item(303, 1).Value2 = 44377
$ws.Cells.Item(303, 2).Value2 = 0
$ws.Cells.Item(303, 3).Value2 = 2
$ws.Cells.Item(303, 4).Value2 = 6.063729800200103
$ws.Cells.Item(304, 1).Value2 = 44378
$ws.Cells.Item(304, 2).Value2 = 0
$ws.Cells.Item(304, 3).Value2 = 1
$ws.Cells.Item(304, 4).Value2 = 3.031864900100051
$ws.Cells.Item(305, 1).Value2 = 44379
$ws.Cells.Item(305, 2).Value2 = 1
$ws.Cells.Item(305, 3).Value2 = 2
$ws.Cells.Item(305, 4).Value2 = 6.063729800200103
$ws.Cells.Item(306, 1).Value2 = 44380
$ws.Cells.Item(306, 2).Value2 = 1
$ws.Cells.Item(306, 3).Value2 = 3
$ws.Cells.Item(306, 4).Value2 = 9.095594700300154
$ws.Cells.Item(307, 1).Value2 = 44381
$ws.Cells.Item(307, 2).Value2 = 0
$ws.Cells.Item(307, 3).Value2 = 3
$ws.Cells.Item(307, 4).Value2 = 9.095594700300154
$ws.Cells.Item(308, 1).Value2 = 44382
$ws.Cells.Item(308, 2).Value2 = 0
$ws.Cells.Item(308, 3).Value2 = 2
$ws.Cells.Item(308, 4).Value2 = 6.063729800200103
$ws.Cells.Item(309, 1).Value2 = 44383
$ws.Cells.Item(309, 2).Value2 = 0
$ws.Cells.Item(309, 3).Value2 = 2
$ws.Cells.Item(309, 4).Value2 = 6.063729800200103
$ws.Cells.Item(310, 1).Value2 = 44384
$ws.Cells.Item(310, 2).Value2 = 0
$ws.Cells.Item(310, 3).Value2 = 2
$ws.Cells.Item(310, 4).Value2 = 6.063729800200103
$ws.Cells.Item(311, 1).Value2 = 44385
$ws.Cells.Item(311, 2).Value2 = 0
$ws.Cells.Item(311, 3).Value2 = 2
$ws.Cells.Item(311, 4).Value2 = 6.063729800200103
$ws.Cells.Item(312, 1).Value2 = 44386
$ws.Cells.Item(312, 2).Value2 = 0
$ws.Cells.Item(312, 3).Value2 = 1
$ws.Cells.Item(312, 4).Value2 = 3.031864900100051
$ws.Cells.Item(313, 1).Value2 = 44387
$ws.Cells.Item(313, 2).Value2 = 0
$ws.Cells.Item(313, 3).Value2 = 0
$ws.Cells.Item(313, 4).Value2 = 0
$ws.Cells.Item(314, 1).Value2 = 44388
$ws.Cells.Item(314, 2).Value2 = 0
$ws.Cells.Item(314, 3).Value2 = 0
$ws.Cells.Item(314, 4).Value2 = 0
$ws.Cells.Item(315, 1).Value2 = 44389
$ws.Cells.Item(315, 2).Value2 = 5
$ws.Cells.Item(315, 3).Value2 = 5
$ws.Cells.Item(315, 4).Value2 = 15.15932450050026
$ws.Cells.Item(316, 1).Value2 = 44390
$ws.Cells.Item(316, 2).Value2 = 0
$ws.Cells.Item(316, 3).Value2 = 5
$ws.Cells.Item(316, 4).Value2 = 15.15932450050026
$ws.Cells.Item(317, 1).Value2 = 44391
$ws.Cells.Item(317, 2).Value2 = 0
$ws.Cells.Item(317, 3).Value2 = 5
$ws.Cells.Item(317, 4).Value2 = 15.15932450050026
$ws.Cells.Item(318, 1).Value2 = 44392
$ws.Cells.Item(318, 2).Value2 = 0
$ws.Cells.Item(318, 3).Value2 = 5
$ws.Cells.Item(318, 4).Value2 = 15.15932450050026
$ws.Cells.Item(319, 1).Value2 = 44393
$ws.Cells.Item(319, 2).Value2 = 0
$ws.Cells.Item(319, 3).Value2 = 5
$ws.Cells.Item(319, 4).Value2 = 15.15932450050026
$ws.Cells.Item(320, 1).Value2 = 44394
$ws.Cells.Item(320, 2).Value2 = 4
$ws.Cells.Item(320, 3).Value2 = 9
$ws.Cells.Item(320, 4).Value2 = 27.28678410090046
$ws.Cells.Item(321, 1).Value2 = 44395
$ws.Cells.Item(321, 2).Value2 = 1
$ws.Cells.Item(321, 3).Value2 = 10
$ws.Cells.Item(321, 4).Value2 = 30.31864900100052
$ws.Cells.Item(322, 1).Value2 = 44396
$ws.Cells.Item(322, 2).Value2 = 0
$ws.Cells.Item(322, 3).Value2 = 5
$ws.Cells.Item(322, 4).Value2 = 15.15932450050026
$ws.Cells.Item(323, 1).Value2 = 44397
$ws.Cells.Item(323, 2).Value2 = 1
$ws.Cells.Item(323, 3).Value2 = 6
$ws.Cells.Item(323, 4).Value2 = 18.19118940060031
$ws.Cells.Item(324, 1).Value2 = 44398
$ws.Cells.Item(324, 2).Value2 = 0
$ws.Cells.Item(324, 3).Value2 = 6
$ws.Cells.Item(324, 4).Value2 = 18.19118940060031
$ws.Cells.Item(325, 1).Value2 = 44399
$ws.Cells.Item(325, 2).Value2 = 5
$ws.Cells.Item(325, 3).Value2 = 11
$ws.Cells.Item(325, 4).Value2 = 33.35051390110057
$ws.Cells.Item(326, 1).Value2 = 44400
$ws.Cells.Item(326, 2).Value2 = 8
$ws.Cells.Item(326, 3).Value2 = 19
$ws.Cells.Item(326, 4).Value2 = 57.60543310190099
$ws.Cells.Item(327, 1).Value2 = 44401
$ws.Cells.Item(327, 2).Value2 = 4
$ws.Cells.Item(327, 3).Value2 = 19
$ws.Cells.Item(327, 4).Value2 = 57.60543310190099
$ws.Cells.Item(328, 1).Value2 = 44402
$ws.Cells.Item(328, 2).Value2 = 5
$ws.Cells.Item(328, 3).Value2 = 23
$ws.Cells.Item(328, 4).Value2 = 69.73289270230119

